$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = "Raichur"
$ws.Range("G11").Value = "Raichur"
$ws.Range("G24").Value = "Raichur"
